# Update countries & provincias Spain
# Applies the daily COVID-19 data refresh to the "Pais" sheet:
#  - Updates the "datos actualizados" timestamp in A1
#  - Updates case counts (columns B-H) for several countries whose
#    numbers moved since the last snapshot
#  - Because rows are kept sorted by total cases (column B) descending,
#    a few countries leapfrogged their neighbours; those rows keep their
#    row number but now hold a different country name + its own stats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = 'Datos actualizados a 28 de Marzo de 2020 a las 20:59'

# --- Straightforward stat refreshes (country keeps its row) ----------
$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 118592
$ws.Range("C4").Value = 14466
$ws.Range("D4").Value = 3224
$ws.Range("E4").Value = 113389
$ws.Range("F4").Value = 2666
$ws.Range("G4").Value = 283
$ws.Range("H4").Value = 1979

$ws.Range("A12").Value = 'Suiza'
$ws.Range("B12").Value = 14076
$ws.Range("C12").Value = 1148
$ws.Range("D12").Value = 1530
$ws.Range("E12").Value = 12282
$ws.Range("F12").Value = 301
$ws.Range("G12").Value = 33
$ws.Range("H12").Value = 264

$ws.Range("A17").Value = 'Turquia'
$ws.Range("B17").Value = 7402
$ws.Range("C17").Value = 1704
$ws.Range("D17").Value = 70
$ws.Range("E17").Value = 7224
$ws.Range("F17").Value = 309
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 108

$ws.Range("A18").Value = 'Canada'
$ws.Range("B18").Value = 5526
$ws.Range("C18").Value = 769
$ws.Range("D18").Value = 354
$ws.Range("E18").Value = 5117
$ws.Range("F18").Value = 120
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 55

$ws.Range("A20").Value = 'Noruega'
$ws.Range("B20").Value = 3998
$ws.Range("C20").Value = 227
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 3969
$ws.Range("F20").Value = 84
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 22

$ws.Range("A64").Value = 'Argelia'
$ws.Range("B64").Value = 454
$ws.Range("C64").Value = 45
$ws.Range("D64").Value = 31
$ws.Range("E64").Value = 394
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 29

$ws.Range("A68").Value = 'Lituania'
$ws.Range("B68").Value = 394
$ws.Range("C68").Value = 36
$ws.Range("D68").Value = 1
$ws.Range("E68").Value = 386
$ws.Range("F68").Value = 2
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 7

$ws.Range("A149").Value = 'Tanzania'
$ws.Range("B149").Value = 14
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 1
$ws.Range("E149").Value = 13
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 0

# --- Pakistan overtakes Rumania (rows 34/35 swap) ---------------------
$ws.Range("A34").Value = 'Pakistan'
$ws.Range("B34").Value = 1495
$ws.Range("C34").Value = 122
$ws.Range("D34").Value = 29
$ws.Range("E34").Value = 1454
$ws.Range("F34").Value = 7
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 12

$ws.Range("A35").Value = 'Rumania'
$ws.Range("B35").Value = 1452
$ws.Range("C35").Value = 160
$ws.Range("D35").Value = 139
$ws.Range("E35").Value = 1283
$ws.Range("F35").Value = 34
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 30

# --- Tunez jumps ahead of Jordania, Republica de Macedonia, Kuwait, ---
# --- Moldavia and Kazajistan (rows 80-85 shift down by one) -----------
$ws.Range("A80").Value = 'Tunez'
$ws.Range("B80").Value = 257
$ws.Range("C80").Value = 30
$ws.Range("D80").Value = 2
$ws.Range("E80").Value = 247
$ws.Range("F80").Value = 10
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 8

$ws.Range("A81").Value = 'Jordania'
$ws.Range("B81").Value = 246
$ws.Range("C81").Value = 11
$ws.Range("D81").Value = 18
$ws.Range("E81").Value = 227
$ws.Range("F81").Value = 3
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 1

$ws.Range("A82").Value = 'Republica de Macedonia'
$ws.Range("B82").Value = 241
$ws.Range("C82").Value = 22
$ws.Range("D82").Value = 3
$ws.Range("E82").Value = 234
$ws.Range("F82").Value = 1
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 4

$ws.Range("A83").Value = 'Kuwait'
$ws.Range("B83").Value = 235
$ws.Range("C83").Value = 10
$ws.Range("D83").Value = 64
$ws.Range("E83").Value = 171
$ws.Range("F83").Value = 11
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 0

$ws.Range("A84").Value = 'Moldavia'
$ws.Range("B84").Value = 231
$ws.Range("C84").Value = 32
$ws.Range("D84").Value = 2
$ws.Range("E84").Value = 227
$ws.Range("F84").Value = 33
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 2

$ws.Range("A85").Value = 'Kazajistan'
$ws.Range("B85").Value = 228
$ws.Range("C85").Value = 78
$ws.Range("D85").Value = 16
$ws.Range("E85").Value = 211
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 1

# --- Guatemala jumps ahead of Isla de Man and Guam (rows 131-133) -----
$ws.Range("A131").Value = 'Guatemala'
$ws.Range("B131").Value = 34
$ws.Range("C131").Value = 6
$ws.Range("D131").Value = 10
$ws.Range("E131").Value = 23
$ws.Range("F131").Value = 1
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 1

$ws.Range("A132").Value = 'Isla de Man'
$ws.Range("B132").Value = 32
$ws.Range("C132").Value = 3
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 32
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 0

$ws.Range("A133").Value = 'Guam'
$ws.Range("B133").Value = 32
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 0
$ws.Range("E133").Value = 31
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 1
